$wb = $excel.ActiveWorkbook

# The "想去人数" (number of people interested) counts were refreshed for the
# generated output. Both the "展览" and "全部类型" sheets contain the same
# rows (3-8) and need the same updated values in column F.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 2524
    $ws.Range("F4").Value = 494
    $ws.Range("F5").Value = 86
    $ws.Range("F6").Value = 6545
    $ws.Range("F7").Value = 385
    $ws.Range("F8").Value = 6
}
